# "udpate apiaas 6 feb" — clear the sample/template values on the
# "TopUp Cashback Nominal" sheet (A2:C2, A3:C3, B18) back to blank while
# keeping them as (empty) text cells, matching the rest of the workbook's
# blank template rows (e.g. A4, E4:E15 already use this convention).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TopUp Cashback Nominal")

$cells = @("A2", "B2", "C2", "A3", "B3", "C3", "B18")

foreach ($addr in $cells) {
    $rng = $ws.Range($addr)
    # A leading apostrophe forces a literal (text) entry so the cell keeps
    # its "string" type even though the content is empty, rather than
    # Excel treating the assignment as "clear the cell".
    $rng.Formula = "'"
    # Drop back to the default style so we don't leave a stray
    # quote-prefix / text-format style on the cell.
    $rng.Style = "Normal"
}
